{"js": "// Find the \"Testing and Evidence Log \" text in the body and color it red\n// (FF0000), matching the target revision where that phrase is highlighted\n// within the sentence \"... will be documented in the Testing and Evidence\n// Log document.\"\nconst body = context.document.body;\nconst searchText = \"Testing and Evidence Log \";\nconst results = body.search(searchText, { matchCase: true, matchWholeWord: false });\nresults.load(\"text\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  // Fallback: some engines trim trailing spaces from search hits, so retry\n  // without the trailing space and recolor just that span.\n  const fallback = body.search(\"Testing and Evidence Log\", { matchCase: true });\n  fallback.load(\"text\");\n  await context.sync();\n  for (const r of fallback.items) {\n    r.font.color = \"#FF0000\";\n  }\n} else {\n  for (const r of results.items) {\n    r.font.color = \"#FF0000\";\n  }\n}\n\nawait context.sync();\n", "ps1": "# Color the phrase \"Testing and Evidence Log \" red (FF0000 / wdColorRed)\n# within the sentence \"... will be documented in the Testing and Evidence\n# Log document.\" to match the target revision.\n\n$d = $word.ActiveDocument\n\n$searchText = \"Testing and Evidence Log \"\n$range = $d.Content\n$found = $range.Find.Execute($searchText)\n\nif (-not $found) {\n    # Fallback in case trailing-space matching behaves differently: retry\n    # without the trailing space.\n    $range = $d.Content\n    $searchText = \"Testing and Evidence Log\"\n    $found = $range.Find.Execute($searchText)\n}\n\nif ($found) {\n    $range.Font.Color = [Microsoft.Office.Interop.Word.WdColor]::wdColorRed\n}\n"}
